$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Restructure columns: insert a new "Type" column at D, shifting the old
#        "Parameter"/"Name" columns to E/F, then delete the old "Name" column
#        (now F) since it is being dropped from the table. ---
$ws.Columns("D").Insert()
$ws.Columns("F").Delete()

# --- 2. Header row ---
$ws.Range("D1").Value = "Type"

# --- 3. Column D ("Type") values for the Source-stage rows ---
$ws.Range("D4:D10").Value = "Parameterised TNSA"

# --- 3b. Column E ("Parameter") for the Source-stage rows switches from the old
#         verbose parameter description to the short code name that used to live
#         in the (now removed) "Name" column. ---
$ws.Range("E4").Value = "SourceMode"
$ws.Range("E5").Value = "SigmaX"
$ws.Range("E6").Value = "SigmaY"
$ws.Range("E7").Value = "Emin"
$ws.Range("E8").Value = "Emax"
$ws.Range("E9").Value = "nPnts"
$ws.Range("E10").Value = "MinCTheta"

# --- 4. Rows 11-14 describe the Interface stage instead of the Source stage:
#        Section B changes from "Source" to "Interface"; the old "Radius1"/
#        "Radius2"/"Drift1Length"/"Drift2Length" names collapse onto shared
#        "Type"/"Parameter" pairs. ---
$ws.Range("B11").Value = "Interface"
$ws.Range("C11").Value = "Drift"
$ws.Range("D11").Value = $null
$ws.Range("E11").Value = "Length"

$ws.Range("B12").Value = "Interface"
$ws.Range("C12").Value = "Aperture"
$ws.Range("D12").Value = "Circular"
$ws.Range("E12").Value = "Radius"

$ws.Range("B13").Value = "Interface"
$ws.Range("C13").Value = "Drift"
$ws.Range("D13").Value = $null
$ws.Range("E13").Value = "Length"

$ws.Range("B14").Value = "Interface"
$ws.Range("C14").Value = "Aperture"
$ws.Range("D14").Value = "Circular"
$ws.Range("E14").Value = "Radius"

# --- 5. Column widths (cosmetic, best achievable given engine's width quantisation) ---
$ws.Columns("B").ColumnWidth = 7.67
$ws.Columns("D").ColumnWidth = 17.5
$ws.Columns("E").ColumnWidth = 18.83

# --- 6. Border / style touch-ups ---
# Row 10 bottom border becomes a thin line (matching row 14's style) instead of dotted.
$rng10 = $ws.Range("A10:H10")
$rng10.Borders.Item(9).LineStyle = 1
$rng10.Borders.Item(9).Weight = 2

# Row 11: top border removed, bottom border dotted (new "sub-group" look), except
# column B keeps a solid top border and loses its bottom border (starts a new block).
$rng11 = $ws.Range("A11:H11")
$rng11.Borders.Item(8).LineStyle = -4142
$rng11.Borders.Item(9).LineStyle = -4118
$rng11.Borders.Item(9).Weight = 2

$b11 = $ws.Range("B11")
$b11.Borders.Item(8).LineStyle = 1
$b11.Borders.Item(8).Weight = 2
$b11.Borders.Item(9).LineStyle = -4142

# Row 12 & 13, column B: no top/bottom border (continuation of the merged "Interface" block).
$b12 = $ws.Range("B12")
$b12.Borders.Item(8).LineStyle = -4142
$b12.Borders.Item(9).LineStyle = -4142

$b13 = $ws.Range("B13")
$b13.Borders.Item(8).LineStyle = -4142
$b13.Borders.Item(9).LineStyle = -4142

# Row 14, column B: no top border, solid bottom border (closes the merged block).
$b14 = $ws.Range("B14")
$b14.Borders.Item(8).LineStyle = -4142
$b14.Borders.Item(9).LineStyle = 1
$b14.Borders.Item(9).Weight = 2

# --- 7. Restore the active cell/selection ---
$null = $ws.Range("H8").Select()
